$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D and G keep their text (string) representation instead of
# being auto-converted to numbers when we assign numeric-looking strings.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "242.89"
$ws.Range("G2").Value = "23"

$ws.Range("D3").Value = "23.08"
$ws.Range("G3").Value = "23"

$ws.Range("G4").Value = "23"

$ws.Range("D5").Value = "0.05916"
$ws.Range("G5").Value = "23"

$ws.Range("G6").Value = "23"

$ws.Range("D7").Value = "6.520"
$ws.Range("G7").Value = "23"

$ws.Range("D8").Value = "0.8095"
$ws.Range("G8").Value = "23"

$ws.Range("D9").Value = "0.9290"
$ws.Range("G9").Value = "23"

$ws.Range("D10").Value = "0.1424"
$ws.Range("G10").Value = "23"

$ws.Range("D11").Value = "0.07437"
$ws.Range("G11").Value = "23"

$ws.Range("D12").Value = "0.03231"
$ws.Range("G12").Value = "23"

$ws.Range("D13").Value = "0.03081"
$ws.Range("G13").Value = "23"

$ws.Range("D14").Value = "0.09355"
$ws.Range("G14").Value = "23"

$ws.Range("D15").Value = "3.873"
$ws.Range("G15").Value = "23"

$ws.Range("D16").Value = "0.001576"
$ws.Range("G16").Value = "23"

$ws.Range("D17").Value = "0.04680"
$ws.Range("G17").Value = "23"

$ws.Range("D18").Value = "0.0005898"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("G18").Value = "23"

$ws.Range("D19").Value = "0.005961"
$ws.Range("G19").Value = "23"

$ws.Range("D20").Value = "0.001255"
$ws.Range("E20").Value = "19BitKanKANBestin24h"
$ws.Range("G20").Value = "23"

$ws.Range("D21").Value = "0.004906"
$ws.Range("G21").Value = "23"

$ws.Range("D22").Value = "0.00006809"
$ws.Range("G22").Value = "23"

$ws.Range("D23").Value = "3.571"
$ws.Range("G23").Value = "23"

$ws.Range("D24").Value = "2.142"
$ws.Range("G24").Value = "23"

$ws.Range("G25").Value = "23"

$ws.Range("G26").Value = "23"

$ws.Range("G27").Value = "23"

$ws.Range("G28").Value = "23"

$ws.Range("G29").Value = "23"

$ws.Range("G30").Value = "23"

$ws.Range("G31").Value = "23"

$ws.Range("G32").Value = "23"

$ws.Range("G33").Value = "23"

$ws.Range("G34").Value = "23"

$ws.Range("G35").Value = "23"

$ws.Range("G36").Value = "23"

$ws.Range("G37").Value = "23"

$ws.Range("G38").Value = "23"

$ws.Range("G39").Value = "23"

$ws.Range("D40").Value = "0.03959"
$ws.Range("G40").Value = "23"

$ws.Range("G41").Value = "23"

$ws.Range("D42").Value = "0.1076"
$ws.Range("G42").Value = "23"

$ws.Range("D43").Value = "0.002574"
$ws.Range("G43").Value = "23"

$ws.Range("D44").Value = "0.008774"
$ws.Range("G44").Value = "23"

$ws.Range("D45").Value = "0.00005252"
$ws.Range("G45").Value = "23"

$ws.Range("G46").Value = "23"

$ws.Range("G47").Value = "23"

$ws.Range("G48").Value = "23"

$ws.Range("G49").Value = "23"

$ws.Range("G50").Value = "23"

$ws.Range("G51").Value = "23"
